$wb = $excel.ActiveWorkbook

# --- hardware sheet: E8 becomes a formula (25/899) instead of the literal 25 ---
$wsHardware = $wb.Worksheets.Item("hardware")
$wsHardware.Range("E8").Formula = "=25/899"

# --- orb_mission sheet: B5 value 1 -> 5, D5 label text updated ---
$wsOrbMission = $wb.Worksheets.Item("orb_mission")
$wsOrbMission.Range("B5").Value = 5
$wsOrbMission.Range("D5").Value = "full orbiter slew maneuver using thrusters"

# --- probe_props sheet: new unit label "kg/m3" in C4 ---
$wsProbeProps = $wb.Worksheets.Item("probe_props")
$wsProbeProps.Range("C4").Value = "kg/m3"

# --- Update the view/selection state on every sheet other than orb_mission first,
#     so that orb_mission ends up as the last-activated (and thus "active") tab. ---
$wsPrimary = $wb.Worksheets.Item("PRIMARY INPUTS")
$wsPrimary.Range("B4").Select()

$wsHardware.Range("E8").Select()

$wsOrbProps = $wb.Worksheets.Item("orb_props")
$wsOrbProps.Range("B8").Select()

$wsProbeProps.Range("C4").Select()

$wsProbeMission = $wb.Worksheets.Item("probe_mission")
$wsProbeMission.Range("A2").Select()

# --- orb_mission becomes the active sheet/tab with B5 selected ---
$wsOrbMission.Activate()
$wsOrbMission.Range("B5").Select()
